# Actualización automática 2025-10-14 09:30:10
#
# ALMEIDA CUATIN JHONATHANN CARLOS registered a new "PORCELANATO" sale of
# 2413.67 against "ORTIZ PEREZ KEVIN DAVID" (row 21 of the detail sheets).
# Ripple this new sale through the dependent monthly / compliance sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "VENTAS POR GRUPO": record the new sale under PORCELANATO (col M) ---
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("M21").Value = 2413.67
# Advisor count footer: one more advisor now has PORCELANATO sales (5 -> 6 of 34)
$wsGrupo.Range("M36").Value = "6 de 34"

# --- Sheet "VENTA MENSUAL": same sale lands in the "octubre" column (col F) ---
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F21").Value = 2413.67
# Column total for octubre grows by the same amount
$wsMensual.Range("F36").Value = 12130.03

# --- Sheet "CUMPLIMIENTO MENSUAL": recompute PORCELANATO compliance row ---
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
# "VENTA" column widened by one character to better fit the updated figures
# (ColumnWidth is in characters; the engine adds a fixed 5/6 padding to the
# stored <col width>, so 13.1666... here lands exactly on width="14")
$wsCumplimiento.Columns.Item(4).ColumnWidth = 13.166666666666666
$wsCumplimiento.Range("D12").Value = 10819.22
$wsCumplimiento.Range("E12").Value = 10882.05
$wsCumplimiento.Range("F12").Value = 0.4985523888694071

# ...and the TOTAL row beneath it
$wsCumplimiento.Range("D14").Value = 12130.03
$wsCumplimiento.Range("E14").Value = 24455.53723718182
$wsCumplimiento.Range("F14").Value = 0.3315523283091885
